$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 7319
$ws.Range("I3").Value = 7490
$ws.Range("J3").Value = 7697
$ws.Range("E4").Value = 2016
$ws.Range("J4").Value = 1673
$ws.Range("J5").Value = 604
$ws.Range("J6").Value = 10511
$ws.Range("E7").Value = 26021
$ws.Range("I7").Value = 26231
$ws.Range("J7").Value = 27804

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J3").Value = 509
$ws.Range("J5").Value = 44
$ws.Range("J6").Value = 646
$ws.Range("J7").Value = 1752

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J2").Value = 165
$ws.Range("J3").Value = 206
$ws.Range("J4").Value = 22
$ws.Range("J6").Value = 149
$ws.Range("J7").Value = 558

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value = 287
$ws.Range("J3").Value = 418
$ws.Range("J6").Value = 449
$ws.Range("J7").Value = 1262

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J6").Value = 263
$ws.Range("J7").Value = 701

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J6").Value = 213
$ws.Range("J8").Value = 1752
$ws.Range("J15").Value = 343
$ws.Range("J18").Value = 225
$ws.Range("J19").Value = 800
$ws.Range("J20").Value = 598
$ws.Range("J25").Value = 143
$ws.Range("J29").Value = 1484
$ws.Range("J31").Value = 289
$ws.Range("J33").Value = 1262
$ws.Range("J36").Value = 375
$ws.Range("J41").Value = 209
$ws.Range("J42").Value = 1190
$ws.Range("J47").Value = 203
$ws.Range("I50").Value = 136
$ws.Range("J50").Value = 165
$ws.Range("J51").Value = 347
$ws.Range("J52").Value = 709
$ws.Range("J54").Value = 550
$ws.Range("J55").Value = 438
$ws.Range("E63").Value = 357
$ws.Range("J63").Value = 83
$ws.Range("J64").Value = 183
$ws.Range("J65").Value = 701
$ws.Range("J67").Value = 1024
$ws.Range("J71").Value = 92
$ws.Range("J77").Value = 196
$ws.Range("J78").Value = 321
$ws.Range("J79").Value = 765
$ws.Range("J83").Value = 558
$ws.Range("J84").Value = 232
$ws.Range("J85").Value = 1136
$ws.Range("J86").Value = 170
$ws.Range("J87").Value = 92
$ws.Range("J88").Value = 295
$ws.Range("J89").Value = 344
$ws.Range("J91").Value = 319
$ws.Range("J94").Value = 309
$ws.Range("J97").Value = 254
$ws.Range("J98").Value = 206
$ws.Range("E101").Value = 26021
$ws.Range("I101").Value = 26231
$ws.Range("J101").Value = 27804

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J3").Value = 70
$ws.Range("J7").Value = 289

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J3").Value = 385
$ws.Range("J7").Value = 1024

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("J2").Value = 69
$ws.Range("J7").Value = 232

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J6").Value = 255
$ws.Range("J7").Value = 550

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 449
$ws.Range("J3").Value = 525
$ws.Range("J6").Value = 375
$ws.Range("J7").Value = 1484

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J2").Value = 196
$ws.Range("J6").Value = 309
$ws.Range("J7").Value = 800

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("J2").Value = 63
$ws.Range("J7").Value = 213

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("J3").Value = 31
$ws.Range("J7").Value = 209

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J2").Value = 247
$ws.Range("J6").Value = 631
$ws.Range("J7").Value = 1190

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J6").Value = 99
$ws.Range("J7").Value = 321

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("J6").Value = 247
$ws.Range("J7").Value = 438

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("J6").Value = 83
$ws.Range("J7").Value = 319

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J2").Value = 217
$ws.Range("J7").Value = 765

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("J3").Value = 46
$ws.Range("J7").Value = 183

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J3").Value = 197
$ws.Range("J6").Value = 175
$ws.Range("J7").Value = 598

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("J6").Value = 105
$ws.Range("J7").Value = 225

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J2").Value = 123
$ws.Range("J3").Value = 121
$ws.Range("J7").Value = 375

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J6").Value = 164
$ws.Range("J7").Value = 309

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("J2").Value = 57
$ws.Range("J6").Value = 29
$ws.Range("J7").Value = 143

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("J2").Value = 46
$ws.Range("J7").Value = 203

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J6").Value = 158
$ws.Range("J7").Value = 343

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("J2").Value = 36
$ws.Range("J7").Value = 206

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("I3").Value = 28
$ws.Range("J3").Value = 43
$ws.Range("I7").Value = 136
$ws.Range("J7").Value = 165

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("J6").Value = 176
$ws.Range("J7").Value = 254

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("J6").Value = 156
$ws.Range("J7").Value = 295

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J3").Value = 99
$ws.Range("J7").Value = 344

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("J4").Value = 91
$ws.Range("J7").Value = 170

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("J4").Value = 32
$ws.Range("J6").Value = 144
$ws.Range("J7").Value = 347

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J5").Value = 28
$ws.Range("J7").Value = 1136

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("J2").Value = 27
$ws.Range("J7").Value = 92

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("J3").Value = 64
$ws.Range("J6").Value = 36
$ws.Range("J7").Value = 196

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J4").Value = 26
$ws.Range("J7").Value = 709

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("J6").Value = 61
$ws.Range("J7").Value = 92
